# Update the "rotation" description text on the INFO sheet (B7), matching
# the corrected/reworded explanation of SUBCALC's rotation convention.
$wb = $excel.ActiveWorkbook
$wsInfo = $wb.Worksheets.Item("INFO")
$wsInfo.Range("B7").Value = "The rotation of the tower in degrees, where 0 is along the positive x axis, increasing counter-clockwise. The SUBCALC program defines zero rotaion along the negative y axis, which is baffling."

# Make the INFO sheet the active/selected sheet again, with B7 selected,
# instead of the template sheet.
$wsInfo.Range("B7").Select()
